$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text so Excel does not auto-convert
# numeric/percentage-looking strings into actual numbers (which would
# lose exact text formatting such as trailing zeros in "332.90" or the
# precise percentage text "0.63%"). Only the specific cells that the
# source data changed are touched, keeping every other cell (including
# header row formatting and untouched rows) exactly as it was.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.63%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.14%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.716"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.78%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08081"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.50%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.030"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.41%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.734"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.75%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.543"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.59%"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.81%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9223"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.88%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1263"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.20%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1944"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.88%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.869"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-7.36%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09411"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.68%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03694"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "5.41%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1053"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.38%"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001307"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.23%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006288"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.66%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.372"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.62%"

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.52%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1421"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.59%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2657"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.95%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04423"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.26%"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.56%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004309"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.23%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001242"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.38%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02872"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "14.70%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05496"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.78%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007778"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.03%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009933"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "10.61%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1418"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.32%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002235"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.21%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01111"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.63%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006835"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.68%"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.66%"

# Row 48
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002287"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "27.50%"

# Row 49
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003010"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-12.98%"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.66%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.66%"
